{"js": "// Update the 25 \"three-digit x one-digit\" multiplication prompts in the\n// single table of the document. The prompts live in the 5 content rows\n// (row indices 0, 4, 9, 14, 19 of the table; the other rows are blank\n// practice rows), 5 cells each, left to right, top to bottom.\n//\n// Each (row, col) cell is addressed positionally (not by searching the\n// whole document for the old text) because several of the old values\n// equal some of the new values used elsewhere in the grid - a global\n// find/replace could match the wrong (already-updated) cell. Scoping\n// the search to the specific cell's body keeps every replacement\n// unambiguous and preserves the existing run/paragraph formatting\n// (font, size, alignment) since only the matched text range is swapped.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// [rowIndex, [[colIndex, oldText, newText], ...]]\nconst rowEdits = [\n  [0, [\n    [0, \"636\u00d72=\", \"375\u00d78=\"],\n    [1, \"706\u00d78=\", \"252\u00d78=\"],\n    [2, \"665\u00d74=\", \"644\u00d76=\"],\n    [3, \"744\u00d73=\", \"845\u00d78=\"],\n    [4, \"612\u00d73=\", \"303\u00d79=\"],\n  ]],\n  [4, [\n    [0, \"686\u00d76=\", \"813\u00d77=\"],\n    [1, \"728\u00d73=\", \"233\u00d76=\"],\n    [2, \"626\u00d73=\", \"853\u00d76=\"],\n    [3, \"907\u00d74=\", \"764\u00d75=\"],\n    [4, \"712\u00d75=\", \"872\u00d74=\"],\n  ]],\n  [9, [\n    [0, \"808\u00d74=\", \"854\u00d77=\"],\n    [1, \"588\u00d76=\", \"920\u00d73=\"],\n    [2, \"252\u00d78=\", \"420\u00d72=\"],\n    [3, \"233\u00d76=\", \"620\u00d74=\"],\n    [4, \"568\u00d76=\", \"915\u00d78=\"],\n  ]],\n  [14, [\n    [0, \"196\u00d72=\", \"530\u00d78=\"],\n    [1, \"590\u00d75=\", \"686\u00d76=\"],\n    [2, \"540\u00d75=\", \"877\u00d74=\"],\n    [3, \"283\u00d78=\", \"341\u00d77=\"],\n    [4, \"591\u00d77=\", \"498\u00d77=\"],\n  ]],\n  [19, [\n    [0, \"913\u00d78=\", \"169\u00d73=\"],\n    [1, \"700\u00d78=\", \"618\u00d72=\"],\n    [2, \"609\u00d74=\", \"841\u00d74=\"],\n    [3, \"461\u00d73=\", \"472\u00d77=\"],\n    [4, \"737\u00d74=\", \"293\u00d76=\"],\n  ]],\n];\n\nconst searchResults = [];\n\nfor (const [rowIndex, cols] of rowEdits) {\n  for (const [colIndex, oldText, newText] of cols) {\n    const cell = table.getCell(rowIndex, colIndex);\n    const results = cell.body.search(oldText, { matchCase: true });\n    results.load(\"items\");\n    searchResults.push({ results, newText });\n  }\n}\n\nawait context.sync();\n\nfor (const { results, newText } of searchResults) {\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Update the 25 \"three-digit x one-digit\" multiplication prompts in the\n# single table of the document. The prompts live in the 5 content rows\n# (rows 1, 5, 10, 15, 20 of the table; the other rows are blank practice\n# rows), 5 cells each, left to right, top to bottom.\n#\n# Each cell is addressed positionally via Table.Cell(row, column) - not\n# by a document-wide Find/Replace of the old text - because several of\n# the old values equal some of the new values used elsewhere in the\n# grid; a global replace could re-match an already-updated cell.\n# Assigning directly to Cell.Range.Text keeps every replacement\n# unambiguous and preserves the existing run/paragraph formatting\n# (font, size, alignment) since Word only swaps the text inside the\n# cell's existing range.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$t.Cell(1, 1).Range.Text = \"375\u00d78=\"\n$t.Cell(1, 2).Range.Text = \"252\u00d78=\"\n$t.Cell(1, 3).Range.Text = \"644\u00d76=\"\n$t.Cell(1, 4).Range.Text = \"845\u00d78=\"\n$t.Cell(1, 5).Range.Text = \"303\u00d79=\"\n\n$t.Cell(5, 1).Range.Text = \"813\u00d77=\"\n$t.Cell(5, 2).Range.Text = \"233\u00d76=\"\n$t.Cell(5, 3).Range.Text = \"853\u00d76=\"\n$t.Cell(5, 4).Range.Text = \"764\u00d75=\"\n$t.Cell(5, 5).Range.Text = \"872\u00d74=\"\n\n$t.Cell(10, 1).Range.Text = \"854\u00d77=\"\n$t.Cell(10, 2).Range.Text = \"920\u00d73=\"\n$t.Cell(10, 3).Range.Text = \"420\u00d72=\"\n$t.Cell(10, 4).Range.Text = \"620\u00d74=\"\n$t.Cell(10, 5).Range.Text = \"915\u00d78=\"\n\n$t.Cell(15, 1).Range.Text = \"530\u00d78=\"\n$t.Cell(15, 2).Range.Text = \"686\u00d76=\"\n$t.Cell(15, 3).Range.Text = \"877\u00d74=\"\n$t.Cell(15, 4).Range.Text = \"341\u00d77=\"\n$t.Cell(15, 5).Range.Text = \"498\u00d77=\"\n\n$t.Cell(20, 1).Range.Text = \"169\u00d73=\"\n$t.Cell(20, 2).Range.Text = \"618\u00d72=\"\n$t.Cell(20, 3).Range.Text = \"841\u00d74=\"\n$t.Cell(20, 4).Range.Text = \"472\u00d77=\"\n$t.Cell(20, 5).Range.Text = \"293\u00d76=\"\n"}
